# Weekly update: insert two new price records (date 2023-10-30 / serial 45229)
# for "Ajo" (Chino / Primera) at the top of the data block (row 275),
# pushing the existing rows 275-365 down to 277-367.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 275 (shifts 275:365 -> 277:367)
$ws.Rows.Item(275).Resize(2).Insert()

# --- New row 275: Ajo / Chino / Primera, $/caja 10 kilos ---
$ws.Range("A275").Value = 9
$ws.Range("B275").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C275").Value = "Metropolitana"
$ws.Range("D275").Value = 45229
$ws.Range("E275").Value = 13
$ws.Range("F275").Value = 100112003
$ws.Range("G275").Value = "Ajo"
$ws.Range("H275").Value = "Chino"
$ws.Range("I275").Value = "Primera"
$ws.Range("J275").Value = 520
$ws.Range("K275").Value = 19000
$ws.Range("L275").Value = 20000
$ws.Range("M275").Value = 19500
$ws.Range("N275").Value = "$/caja 10 kilos"
$ws.Range("O275").Value = "China"
$ws.Range("P275").Value = 1950
$ws.Range("Q275").Value = 10
$ws.Range("R275").Value = "Hortaliza"

# --- New row 276: Ajo / Chino / Primera, $/malla 10 kilos ---
$ws.Range("A276").Value = 9
$ws.Range("B276").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C276").Value = "Metropolitana"
$ws.Range("D276").Value = 45229
$ws.Range("E276").Value = 13
$ws.Range("F276").Value = 100112003
$ws.Range("G276").Value = "Ajo"
$ws.Range("H276").Value = "Chino"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 340
$ws.Range("K276").Value = 18000
$ws.Range("L276").Value = 19000
$ws.Range("M276").Value = 18500
$ws.Range("N276").Value = "$/malla 10 kilos"
$ws.Range("O276").Value = "China"
$ws.Range("P276").Value = 1850
$ws.Range("Q276").Value = 10
$ws.Range("R276").Value = "Hortaliza"
